$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "Any TCV reported (not specified how defined in codebook)"
# (old row 5) entirely — this shifts rows 6-8 up to 5-7 and drops the now-unused
# shared string along with it.
$ws.Rows.Item(5).Delete()

# Update selection to match the new target state
$ws.Range("F10").Select()
